# "Generate Report for Handback" - refresh the generated timestamps that
# are written into the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date for 54fdd9d3-...-1809dc437207.md
# (shared between Overview!G3 and de-de!H3)
$wsOverview.Range("G3").Value = "2016-08-15 14:44:11"
$wsDeDe.Range("H3").Value = "2016-08-15 14:44:11"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H3").Value = "2016-08-15 14:44:00"
$wsZhCn.Range("K3").Value = "2016-08-15 14:44:30"

# de-de: Correspond Handback DateTime
$wsDeDe.Range("K3").Value = "2016-08-15 14:44:37"
